$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells are written in this precise order so the shared-string table comes
# out in the same sequence the source workbook uses.
$ws.Range("A1").Value = "Name"
$ws.Range("A2").Value = "Veronica Tjan"

$ws.Range("B1").Value = "Age"

$ws.Range("C1").Value = "Gender"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "F"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "20"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "veronicatjan@hotmail.com"

$ws.Range("E1").Value = "Phone"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "13250808969"

$ws.Range("A3").Value = "Sarah Huang"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "23"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "F"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "sarahhuang@gmail.com"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "13235812152"

$ws.Range("D1").Value = "Email"

# Hyperlinks on the email cells
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:veronicatjan@hotmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:sarahhuang@gmail.com")

# Leftover formatted (empty) cell at E9
$ws.Range("E9").HorizontalAlignment = -4131
$ws.Range("E9").Value = ""

# Column widths (Excel stores column width quantised to whole pixels, so the
# ColumnWidth we assign is chosen to round-trip to the saved "characters"
# width closest to the target: col A -> 14.6640625, col D -> 23, col E -> 17.5)
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(4).ColumnWidth = 22.285714285714285
$ws.Columns.Item(5).ColumnWidth = 16.857142857142854

# Selection matches the saved selection state
$ws.Range("F9").Select()
